$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D values are set via NumberFormat "@" + ClearFormats() so that numeric-looking
# strings (e.g. "1.011") are stored as text, matching the original inline-string cells,
# without leaving a residual text style behind.
$priceUpdates = @{
    2 = '27.196.45'
    3 = '1.828.17'
    4 = '1.011'
    5 = '313.77'
    7 = '0.4720'
    8 = '0.3656'
    9 = '0.07419'
    10 = '0.8827'
    11 = '20.43'
    12 = '1.893.71'
    13 = '0.07315'
    14 = '5.403'
    15 = '93.46'
    16 = '6.526'
    18 = '0.000008775'
    20 = '27.673.72'
    21 = '14.70'
    22 = '5.269'
    23 = '10.62'
    24 = '2.117.52'
    25 = '1.887'
    26 = '151.72'
    27 = '18.58'
    28 = '2.148'
    29 = '5.199'
    30 = '116.40'
    31 = '0.08958'
    32 = '1.169'
    33 = '0.7434'
    34 = '4.520'
    35 = '2.947'
    37 = '1.091'
    38 = '0.05320'
    40 = '2.422'
    41 = '2.945'
    42 = '7.206'
    43 = '0.5276'
    44 = '0.1649'
    45 = '8.444'
    46 = '0.4895'
    47 = '10.45'
    49 = '104.97'
    50 = '1.654'
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
    $cell.ClearFormats()
}

$volumeUpdates = @{
    2 = '  +0.88%  '
    3 = '  +0.54%  '
    4 = '  +0.81%  '
    5 = '  +1.23%  '
    6 = '  +0.76%  '
    7 = '  +1.15%  '
    8 = '  -0.20%  '
    9 = '  +0.84%  '
    10 = '  +1.19%  '
    11 = '  +0.81%  '
    12 = '  +3.69%  '
    13 = '  +2.84%  '
    14 = '  -0.15%  '
    15 = '  +2.10%  '
    16 = '  +0.21%  '
    17 = '  +0.47%  '
    18 = '  +0.68%  '
    19 = '  +0.80%  '
    20 = '  +2.58%  '
    21 = '  +0.43%  '
    22 = '  -0.49%  '
    23 = '  +0.24%  '
    24 = '  +3.26%  '
    25 = '  -0.35%  '
    26 = '  +0.58%  '
    27 = '  +0.79%  '
    28 = '  -0.28%  '
    29 = '  -1.00%  '
    30 = '  -0.70%  '
    31 = '  +0.66%  '
    32 = '  +0.60%  '
    33 = '  -2.04%  '
    34 = '  +0.47%  '
    35 = '  +1.18%  '
    36 = '  +0.87%  '
    37 = '  -0.08%  '
    38 = '  +0.38%  '
    39 = '  +0.15%  '
    40 = '  +2.04%  '
    41 = '  -0.90%  '
    42 = '  +0.08%  '
    43 = '  -0.43%  '
    44 = '  -0.12%  '
    45 = '  -0.01%  '
    46 = '  +0.44%  '
    47 = '  -0.55%  '
    48 = '  +0.86%  '
    49 = '  +1.38%  '
    50 = '  -0.70%  '
    51 = '  +0.04%  '
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}
